$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MPD"
$ws.Range("B2").Value = "Makanan Pedas"
$ws.Range("A3").Value = "MPM"
$ws.Range("B3").Value = "Makanan Pedas Manis"
